# Teacher To-Do Mini Intro: insert 5 new content slides and reorder the deck.
$p = $ppt.ActivePresentation

# --- 1. "Introduction" slide, inserted right after "Team Introduction" (position 3) ---
$sIntro = $p.Slides.Add(3, 2)

$introTitle = $sIntro.Shapes.Item(1)
$introTitle.TextFrame.TextRange.Text = "Introduction"

$introBody = $sIntro.Shapes.Item(2)
$introBody.Left   = 565148 / 12700
$introBody.Top    = 2174052 / 12700
$introBody.Width  = 9959077 / 12700
$introBody.Height = 3967955 / 12700

$tr = $introBody.TextFrame.TextRange
$tr.Text = "The education of children is critical to the future success of our entire society. In Arizona alone, there are an estimated 2,600 teaching positions that remain open this school year. In order "
$tr.Font.Size = 20
$tr.Font.Bold = $false

$run2 = $tr.InsertAfter("to help fill this gap and ensure the success of Arizona students, the Arizona Department of Education has started a new program that allows undergraduate education students the opportunity to fill some of these vacant positions and gain hands on experience in the classroom while finishing their degree.")
$run2.Font.Size = 20

$para2 = $introBody.TextFrame.TextRange.Paragraphs(2)
$para2.Text = "`rThe College of Education is tasked with assisting NAU students who would like to enter these vacant positions. This involves ensuring that students meet a variety of requirements set out by AZDE, gathering supporting documentation and ensuring that each student has the support they need to be successful.`r"
$paras = $introBody.TextFrame.TextRange.Paragraphs()
$paras.Item(2).Font.Size = 20
$paras.Item(2).Font.Bold = $false

# --- 2. "Client Intro" slide, inserted right after "Introduction" (position 4) ---
$sClient = $p.Slides.Add(4, 2)
$sClient.Shapes.Item(1).TextFrame.TextRange.Text = "Client Intro"

# --- 3. "Proposed Solution" slide, inserted after "Current Solution Overview/Problems" (position 6) ---
$sProposed = $p.Slides.Add(6, 2)
$sProposed.Shapes.Item(1).TextFrame.TextRange.Text = "Proposed Solution"

# --- 4. "Conclusion" slide, inserted right after "Proposed Solution" (position 7) ---
$sConclusion = $p.Slides.Add(7, 2)
$sConclusion.Shapes.Item(1).TextFrame.TextRange.Text = "Conclusion"

# --- 5. "Plan for Development" slide, inserted before "Conclusion" (position 7) ---
$sPlan = $p.Slides.Add(7, 2)
$sPlan.Shapes.Item(1).TextFrame.TextRange.Text = "Plan for Development"

$planBody = $sPlan.Shapes.Item(2)
$planBody.Left   = 565149 / 12700
$planBody.Top    = 2211536 / 12700
$planBody.Width  = 10597432 / 12700
$planBody.Height = 3441743 / 12700

$planTr = $planBody.TextFrame.TextRange
$planTr.Text = "We plan to have bi-weekly meetings with our client to have high-level discussions or requirements and obtain feedback on design choices.`rWe have already begun our technical investigation and have started acquiring some resources that may take longer to obtain like CAS integration.`rBecause we are handling student records, we must take student privacy into account to ensure data integrity."
